$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet tab to reflect the new "through" date
$ws.Name = "Through 2022-07-22"

# Update the header label cell (B1) to reflect the new "through" date
$ws.Range("B1").Value = "July 2022 (through July 22)"

# Update existing counts (incremented by 1 for the new day's data)
$ws.Range("B2").Value = 13
$ws.Range("P2").Value = 6
$ws.Range("AD2").Value = 10
$ws.Range("P3").Value = 5
$ws.Range("AK4").Value = 3
$ws.Range("AR6").Value = 4
$ws.Range("I8").Value = 7
$ws.Range("P8").Value = 16
$ws.Range("B20").Value = 5
$ws.Range("AK20").Value = 3
$ws.Range("B29").Value = 7
$ws.Range("AK44").Value = 2
$ws.Range("P45").Value = 2
$ws.Range("W45").Value = 2
$ws.Range("AD47").Value = 4
$ws.Range("AK47").Value = 2
$ws.Range("B52").Value = 4
$ws.Range("I78").Value = 3
$ws.Range("AD94").Value = 2

# Add new counts for cells that previously had no data
$ws.Range("I37").Value = 1
$ws.Range("P38").Value = 1
$ws.Range("AR39").Value = 1
$ws.Range("B45").Value = 2
$ws.Range("P49").Value = 1
$ws.Range("B56").Value = 1
$ws.Range("AD83").Value = 1
$ws.Range("I89").Value = 1
